$d = $word.ActiveDocument

$d.Content.Find.Execute("18÷6=3, 0", $true, $false, $false, $false, $false, $true, 1, $false, "56÷3=18, 2", 2)
$d.Content.Find.Execute("25÷7=3, 4", $true, $false, $false, $false, $false, $true, 1, $false, "41÷6=6, 5", 2)
$d.Content.Find.Execute("17÷2=8, 1", $true, $false, $false, $false, $false, $true, 1, $false, "49÷5=9, 4", 2)
$d.Content.Find.Execute("30÷3=10, 0", $true, $false, $false, $false, $false, $true, 1, $false, "85÷4=21, 1", 2)
$d.Content.Find.Execute("37÷7=5, 2", $true, $false, $false, $false, $false, $true, 1, $false, "18÷7=2, 4", 2)
$d.Content.Find.Execute("58÷2=29, 0", $true, $false, $false, $false, $false, $true, 1, $false, "25÷8=3, 1", 2)
$d.Content.Find.Execute("17÷8=2, 1", $true, $false, $false, $false, $false, $true, 1, $false, "92÷8=11, 4", 2)
$d.Content.Find.Execute("13÷4=3, 1", $true, $false, $false, $false, $false, $true, 1, $false, "21÷6=3, 3", 2)
$d.Content.Find.Execute("42÷4=10, 2", $true, $false, $false, $false, $false, $true, 1, $false, "12÷3=4, 0", 2)
$d.Content.Find.Execute("91÷2=45, 1", $true, $false, $false, $false, $false, $true, 1, $false, "61÷8=7, 5", 2)
$d.Content.Find.Execute("72÷3=24, 0", $true, $false, $false, $false, $false, $true, 1, $false, "22÷4=5, 2", 2)
$d.Content.Find.Execute("23÷9=2, 5", $true, $false, $false, $false, $false, $true, 1, $false, "93÷4=23, 1", 2)
$d.Content.Find.Execute("90÷2=45, 0", $true, $false, $false, $false, $false, $true, 1, $false, "48÷7=6, 6", 2)
$d.Content.Find.Execute("24÷6=4, 0", $true, $false, $false, $false, $false, $true, 1, $false, "80÷5=16, 0", 2)
$d.Content.Find.Execute("98÷8=12, 2", $true, $false, $false, $false, $false, $true, 1, $false, "54÷8=6, 6", 2)
$d.Content.Find.Execute("68÷4=17, 0", $true, $false, $false, $false, $false, $true, 1, $false, "55÷7=7, 6", 2)
$d.Content.Find.Execute("47÷4=11, 3", $true, $false, $false, $false, $false, $true, 1, $false, "46÷3=15, 1", 2)
$d.Content.Find.Execute("46÷2=23, 0", $true, $false, $false, $false, $false, $true, 1, $false, "54÷7=7, 5", 2)
$d.Content.Find.Execute("35÷4=8, 3", $true, $false, $false, $false, $false, $true, 1, $false, "56÷9=6, 2", 2)
$d.Content.Find.Execute("61÷5=12, 1", $true, $false, $false, $false, $false, $true, 1, $false, "94÷8=11, 6", 2)
$d.Content.Find.Execute("71÷9=7, 8", $true, $false, $false, $false, $false, $true, 1, $false, "67÷7=9, 4", 2)
$d.Content.Find.Execute("21÷4=5, 1", $true, $false, $false, $false, $false, $true, 1, $false, "34÷3=11, 1", 2)
$d.Content.Find.Execute("59÷3=19, 2", $true, $false, $false, $false, $false, $true, 1, $false, "96÷7=13, 5", 2)
$d.Content.Find.Execute("12÷4=3, 0", $true, $false, $false, $false, $false, $true, 1, $false, "16÷6=2, 4", 2)
$d.Content.Find.Execute("95÷5=19, 0", $true, $false, $false, $false, $false, $true, 1, $false, "90÷5=18, 0", 2)
